$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("systems")

$ws.Range("A7").Value = "teste"
$ws.Range("B7").Value = "Teste"

$ws.Range("A8").Value = "teste2"
$ws.Range("B8").Value = "Teste2"
